# Daily attendance processing - reverse the order of names/emails listed
# in the "Recorded By" column (G) for every data row on the active sheet.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = 7
    $current = $cell.Value2

    if ($current -ne $null -and $current -ne "") {
        $parts = $current -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newValue = $reversed -join ", "
            if ($newValue -ne $current) {
                $cell.Value2 = $newValue
            }
        }
    }
}
